$d = $word.ActiveDocument

$replacements = @(
    @("2024-07-20 Saturday", "2024-07-21 Sunday"),
    @("11×41=451", "11×40=440"),
    @("67×11=737", "26×82=2132"),
    @("60×81=4860", "16×65=1040"),
    @("61×75=4575", "29×70=2030"),
    @("99×99=9801", "84×35=2940"),
    @("26×96=2496", "66×77=5082"),
    @("82×45=3690", "11×97=1067"),
    @("43×82=3526", "43×91=3913"),
    @("18×56=1008", "52×77=4004"),
    @("16×73=1168", "81×21=1701"),
    @("36×13=468", "63×16=1008"),
    @("96×28=2688", "83×44=3652"),
    @("46×90=4140", "67×71=4757"),
    @("34×61=2074", "44×20=880"),
    @("40×26=1040", "40×91=3640"),
    @("91×99=9009", "52×82=4264"),
    @("97×33=3201", "51×15=765"),
    @("25×53=1325", "86×19=1634"),
    @("20×63=1260", "61×94=5734"),
    @("15×14=210", "85×72=6120"),
    @("82×59=4838", "38×61=2318"),
    @("28×25=700", "80×17=1360"),
    @("15×37=555", "65×13=845"),
    @("29×14=406", "89×25=2225"),
    @("15×11=165", "96×69=6624")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
